$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "Deposit"
$ws.Range("P3").Value = "Roobic"
$ws.Range("T2").Value = 100
$ws.Range("T3").Value = 50

$ws.Range("T3").Select()
